$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C18").Value = 8217
$ws.Range("C19:C22").Value = 8149
$ws.Range("C23:C49").Value = 8085
$ws.Range("C50:C109").Value = 7590
$ws.Range("C110:C148").Value = 7573
$ws.Range("C149:C252").Value = 7569
